$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.060.88"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.155.12"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.87"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.04"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.153.86"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.446"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("E12").Value = "  +4.25%  "
$ws.Range("D13").Value = "3.689.10"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.79"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "58.095.94"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "3.146.45"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.01"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "358.39"
$ws.Range("E22").Value = "  +6.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.96"
$ws.Range("E24").Value = "  +3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.513"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.42"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.33"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("E35").Value = "  +6.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.96"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.22"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.25"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.29"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0674"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.64"
$ws.Range("E41").Value = "  +10.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.13"
$ws.Range("E42").Value = "  +6.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.708"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").Value = "3.188.85"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0275"
$ws.Range("E45").Value = "  +6.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.79"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "2.337.02"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.10"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.50"
$ws.Range("E51").Value = "  -1.72%  "
